$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.01581778927454889, 0.1769645570152142, 0.06154514893109243, 0.2480829476830127, 0.2569240011290105)
    3  = @(0.02195526675411769, 0.2163586422281868, 0.08385254728376079, 0.2895730430888911, 0.3005296263947633)
    4  = @(-0.0193087473139601, 0.2049953786970008, 0.1019947020594539, 0.3193660940980647, 0.3329567938172152)
    5  = @(-0.07226700651545802, 0.1675871387573901, 0.06404507904842406, 0.2530712924225584, 0.2543714109319102)
    6  = @(-0.1025096706313666, 0.165819033310199, 0.06384594992382764, 0.252677561179911, 0.2434422526798966)
    7  = @(-0.1148643701330149, 0.1835058106740494, 0.06418979330833566, 0.2533570470863908, 0.2395213268272089)
    8  = @(-0.09270844093133584, 0.2111731646910383, 0.07863255137221827, 0.2804149628179963, 0.2899055632835616)
    9  = @(-0.2246702250377167, 0.2246702250377167, 0.1088070602894433, 0.3298591521990004, 0.2957964256146743)
    10 = @(-0.09344167408882031, 0.09344167408882031, 0.008731346456521313, 0.09344167408882031, $null)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    if ($null -ne $values[4]) {
        $ws.Range("F$row").Value = $values[4]
    }
}
